$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume(1h) data (refreshed crypto feed snapshot)
$ws.Range("D2").Value = '24.889.30'
$ws.Range("E2").Value = '  +1.76%  '
$ws.Range("D3").Value = '1.669.17'
$ws.Range("E3").Value = '  +0.72%  '
$ws.Range("D4").Value = '''1.003'
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '''330.38'
$ws.Range("E5").Value = '  +7.29%  '
$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '  +0.23%  '
$ws.Range("D7").Value = '''0.3647'
$ws.Range("E7").Value = '  +0.48%  '
$ws.Range("D8").Value = '''46.65'
$ws.Range("E8").Value = '  -1.37%  '
$ws.Range("D9").Value = '''0.3225'
$ws.Range("E9").Value = '  -1.36%  '
$ws.Range("D10").Value = '''1.140'
$ws.Range("E10").Value = '  +1.20%  '
$ws.Range("D11").Value = '''0.07054'
$ws.Range("E11").Value = '  +1.30%  '
$ws.Range("D12").Value = '''1.001'
$ws.Range("E12").Value = '  +0.29%  '
$ws.Range("D13").Value = '''6.067'
$ws.Range("E13").Value = '  +2.14%  '
$ws.Range("D14").Value = '''19.57'
$ws.Range("E14").Value = '  +1.20%  '
$ws.Range("D15").Value = '1.662.87'
$ws.Range("E15").Value = '  +0.82%  '
$ws.Range("D16").Value = '''6.602'
$ws.Range("E16").Value = '  -0.28%  '
$ws.Range("D17").Value = '''0.00001045'
$ws.Range("E17").Value = '  +0.19%  '
$ws.Range("D18").Value = '''0.06543'
$ws.Range("E18").Value = '  +0.34%  '
$ws.Range("E19").Value = '  +0.33%  '
$ws.Range("D20").Value = '''78.56'
$ws.Range("E20").Value = '  +2.83%  '
$ws.Range("D21").Value = '''15.81'
$ws.Range("E21").Value = '  +0.38%  '
$ws.Range("D22").Value = '''5.909'
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("D23").Value = '''12.89'
$ws.Range("E23").Value = '  +1.95%  '
$ws.Range("D24").Value = '24.879.44'
$ws.Range("E24").Value = '  +1.89%  '
$ws.Range("D25").Value = '''2.443'
$ws.Range("E25").Value = '  -0.55%  '
$ws.Range("D26").Value = '''2.396'
$ws.Range("E26").Value = '  +2.80%  '
$ws.Range("D27").Value = '''148.06'
$ws.Range("E27").Value = '  +1.34%  '
$ws.Range("D28").Value = '''18.67'
$ws.Range("E28").Value = '  +1.46%  '
$ws.Range("D29").Value = '1.847.61'
$ws.Range("E29").Value = '  +0.71%  '
$ws.Range("D30").Value = '''125.36'
$ws.Range("E30").Value = '  +1.04%  '
$ws.Range("D31").Value = '''1.174'
$ws.Range("E31").Value = '  -2.34%  '
$ws.Range("E32").Value = '  +0.63%  '
$ws.Range("D33").Value = '''5.725'
$ws.Range("E33").Value = '  +2.27%  '
$ws.Range("D34").Value = '''0.08450'
$ws.Range("E34").Value = '  +1.35%  '
$ws.Range("E35").Value = '  -2.35%  '
$ws.Range("D36").Value = '''12.23'
$ws.Range("E36").Value = '  -1.31%  '
$ws.Range("D37").Value = '''5.138'
$ws.Range("E37").Value = '  -1.33%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '''0.06027'
$ws.Range("E38").Value = '  -0.57%  '
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = '''1.228'
$ws.Range("E39").Value = '  +1.95%  '
$ws.Range("D40").Value = '''0.02232'
$ws.Range("E40").Value = '  +1.58%  '
$ws.Range("D41").Value = '''0.2082'
$ws.Range("E41").Value = '  +1.62%  '
$ws.Range("D42").Value = '''8.227'
$ws.Range("E42").Value = '  +0.55%  '
$ws.Range("E43").Value = '  +0.17%  '
$ws.Range("D44").Value = '''0.5935'
$ws.Range("E44").Value = '  +1.40%  '
$ws.Range("D45").Value = '''13.58'
$ws.Range("E45").Value = '  +7.90%  '
$ws.Range("D46").Value = '''3.851'
$ws.Range("E46").Value = '  +3.15%  '
$ws.Range("D47").Value = '''0.5719'
$ws.Range("E47").Value = '  +2.54%  '
$ws.Range("D48").Value = '''124.94'
$ws.Range("E48").Value = '  +2.34%  '
$ws.Range("D49").Value = '''1.958'
$ws.Range("E49").Value = '  +1.02%  '
$ws.Range("D50").Value = '''0.06983'
$ws.Range("E50").Value = '  +1.10%  '
$ws.Range("D51").Value = '''1.186'
$ws.Range("E51").Value = '  +2.90%  '
